$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 202
$ws1.Range("F5").Value = 247
$ws1.Range("F7").Value = 75
$ws1.Range("F8").Value = 261
$ws1.Range("F11").Value = 34
$ws1.Range("F12").Value = 102
$ws1.Range("F13").Value = 2314
$ws1.Range("F17").Value = 532
$ws1.Range("F22").Value = 1796
$ws1.Range("F23").Value = 3922
$ws1.Range("F26").Value = 1167
$ws1.Range("F28").Value = 2074
$ws1.Range("F32").Value = 100
$ws1.Range("F33").Value = 287
$ws1.Range("F36").Value = 684

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 27

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 202
$ws4.Range("F5").Value = 247
$ws4.Range("F7").Value = 75
$ws4.Range("F8").Value = 261
$ws4.Range("F11").Value = 34
$ws4.Range("F12").Value = 102
$ws4.Range("F13").Value = 2314
$ws4.Range("F15").Value = 27
$ws4.Range("F18").Value = 532
$ws4.Range("F23").Value = 1797
$ws4.Range("F24").Value = 3922
$ws4.Range("F27").Value = 1167
$ws4.Range("F29").Value = 2074
$ws4.Range("F33").Value = 100
$ws4.Range("F34").Value = 287
$ws4.Range("F37").Value = 684
